$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the country names shown on rows 209 and 210 -----------------
# (Groenlandia / Islas Malvinas traded places in the author's source data;
#  the numeric stats on those two rows are identical so only the labels move.)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Islas Malvinas"

# --- Refresh the "last updated" timestamp caption ----------------------
$ws.Range("A1").Value = "Datos actualizados a 6 de Julio de 2020 a las 23:32"

# --- Row 4 (Estados Unidos) --------------------------------------------
$ws.Range("B4").Value = 3022595
$ws.Range("C4").Value = 39667
$ws.Range("D4").Value = 1307683
$ws.Range("E4").Value = 1582119
$ws.Range("G4").Value = 224
$ws.Range("H4").Value = 132793

# --- Row 8 (Peru) --------------------------------------------------------
$ws.Range("B8").Value = 305703
$ws.Range("C8").Value = 2985
$ws.Range("D8").Value = 197619
$ws.Range("E8").Value = 97312
$ws.Range("G8").Value = 183
$ws.Range("H8").Value = 10772

# --- Row 18 (Alemania) ---------------------------------------------------
$ws.Range("B18").Value = 198057
$ws.Range("C18").Value = 499
$ws.Range("E18").Value = 6765
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 9092

# --- Row 28 ---------------------------------------------------------------
$ws.Range("B28").Value = 73061
$ws.Range("C28").Value = 44
$ws.Range("G28").Value = 13
$ws.Range("H28").Value = 5433

# --- Row 50 -----------------------------------------------------------------
$ws.Range("B50").Value = 29821
$ws.Range("C50").Value = 454
$ws.Range("D50").Value = 25178
$ws.Range("E50").Value = 4545

# --- Row 70 -----------------------------------------------------------------
$ws.Range("B70").Value = 10966
$ws.Range("C70").Value = 194
$ws.Range("D70").Value = 5384
$ws.Range("E70").Value = 5507
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 75

# --- Row 104 -----------------------------------------------------------------
$ws.Range("B104").Value = 3006
$ws.Range("C104").Value = 9
$ws.Range("D104").Value = 1051
$ws.Range("E104").Value = 1863

# --- Row 191 -----------------------------------------------------------------
$ws.Range("B191").Value = 61
$ws.Range("C191").Value = 4
$ws.Range("E191").Value = 31
$ws.Range("G191").Value = 1
$ws.Range("H191").Value = 3
